# Daily attendance processing - 2025-10-16 16:52:42
# Normalize the "Recorded By" (column G) values on the active sheet by
# sorting the comma-separated list of recorders alphabetically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "" -and $val -like "*,*") {
        $parts = $val -split ", "
        $sortedParts = $parts | Sort-Object
        $newVal = $sortedParts -join ", "

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

Write-Output "Recorded By column sorted."
